$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": append a new date column CN (13-sep) with hourly prices
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (CM1) onto the new header
# cell (CN1) so the new column keeps the bold/centered/bordered header style,
# then set its text.
$wsPrix.Range("CM1:CM1").Copy($wsPrix.Cells.Item(1, 92))
$wsPrix.Cells.Item(1, 92).Value = "13-sep"

$cnValues = @{
    2  = 9.449999999999999
    3  = 9.17
    4  = 10.35
    5  = 7.39
    6  = 7
    7  = 7.75
    8  = 10.55
    9  = 9.630000000000001
    10 = 13.63
    11 = 12.42
    12 = 9.619999999999999
    13 = 2.37
    14 = 0
    15 = -0.01
    16 = -0.01
    17 = 0
    18 = 0
    19 = 0.65
    20 = 14
    21 = 17.68
    22 = 25.27
    23 = 10.51
    24 = 29.45
    25 = 25.73
}

foreach ($row in 2..25) {
    $wsPrix.Cells.Item($row, 92).Value = $cnValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append one new row (89) with the latest daily price
# ---------------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Cells.Item(89, 1).NumberFormat = "@"
$wsGaz.Cells.Item(89, 1).Value = "2025-09-11"
$wsGaz.Cells.Item(89, 1).Style = "Normal"
$wsGaz.Cells.Item(89, 2).Value = 31.8

# ---------------------------------------------------------------------------
# Sheet "CO2": append one new row (89) with the latest daily price
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Cells.Item(89, 1).NumberFormat = "@"
$wsCo2.Cells.Item(89, 1).Value = "2025-09-11"
$wsCo2.Cells.Item(89, 1).Style = "Normal"
$wsCo2.Cells.Item(89, 2).Value = 75.25
